# Update the SharePoint/OneDrive customXml metadata parts (item2.xml schema
# definitions, item3.xml document-management values and the itemProps2.xml
# datastore descriptor) to the refreshed content-type / field schema that
# ships with the latest XML sample upload. This mirrors what PowerPoint does
# when it resynchronises the "Document Management" custom XML parts attached
# to the file (new field order/indexes, new list GUIDs, new content type id).

$p = $ppt.ActivePresentation
$parts = $p.CustomXMLParts

$item2Xml = '<ct:contentTypeSchema xmlns:ct="http://schemas.microsoft.com/office/2006/metadata/contentType" xmlns:ma="http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes" ct:_="" ma:_="" ma:contentTypeName="Document" ma:contentTypeID="0x010100E62A33F58702334EBF7311526B3DE3F0" ma:contentTypeVersion="15" ma:contentTypeDescription="Create a new document." ma:contentTypeScope="" ma:versionID="977969c0047b5e5be8638840a0ae1307"><xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:ns2="25558cfc-09cb-4f21-a339-24113151f02d" xmlns:ns3="67e6609d-4da6-43f2-b73c-9b749e56ecda" targetNamespace="http://schemas.microsoft.com/office/2006/metadata/properties" ma:root="true" ma:fieldsID="6efdbf2bade035783d617b88be4656c9" ns2:_="" ns3:_=""><xsd:import namespace="25558cfc-09cb-4f21-a339-24113151f02d"/><xsd:import namespace="67e6609d-4da6-43f2-b73c-9b749e56ecda"/><xsd:element name="properties"><xsd:complexType><xsd:sequence><xsd:element name="documentManagement"><xsd:complexType><xsd:all><xsd:element ref="ns2:MediaServiceMetadata" minOccurs="0"/><xsd:element ref="ns2:MediaServiceFastMetadata" minOccurs="0"/><xsd:element ref="ns2:MediaServiceAutoKeyPoints" minOccurs="0"/><xsd:element ref="ns2:MediaServiceKeyPoints" minOccurs="0"/><xsd:element ref="ns3:SharedWithUsers" minOccurs="0"/><xsd:element ref="ns3:SharedWithDetails" minOccurs="0"/><xsd:element ref="ns2:MediaServiceAutoTags" minOccurs="0"/><xsd:element ref="ns2:MediaServiceOCR" minOccurs="0"/><xsd:element ref="ns2:MediaServiceGenerationTime" minOccurs="0"/><xsd:element ref="ns2:MediaServiceEventHashCode" minOccurs="0"/><xsd:element ref="ns2:MediaServiceDateTaken" minOccurs="0"/><xsd:element ref="ns2:MediaLengthInSeconds" minOccurs="0"/><xsd:element ref="ns2:lcf76f155ced4ddcb4097134ff3c332f" minOccurs="0"/><xsd:element ref="ns3:TaxCatchAll" minOccurs="0"/></xsd:all></xsd:complexType></xsd:element></xsd:sequence></xsd:complexType></xsd:element></xsd:schema><xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="25558cfc-09cb-4f21-a339-24113151f02d" elementFormDefault="qualified"><xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/><xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/><xsd:element name="MediaServiceMetadata" ma:index="8" nillable="true" ma:displayName="MediaServiceMetadata" ma:hidden="true" ma:internalName="MediaServiceMetadata" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Note"/></xsd:simpleType></xsd:element><xsd:element name="MediaServiceFastMetadata" ma:index="9" nillable="true" ma:displayName="MediaServiceFastMetadata" ma:hidden="true" ma:internalName="MediaServiceFastMetadata" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Note"/></xsd:simpleType></xsd:element><xsd:element name="MediaServiceAutoKeyPoints" ma:index="10" nillable="true" ma:displayName="MediaServiceAutoKeyPoints" ma:hidden="true" ma:internalName="MediaServiceAutoKeyPoints" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Note"/></xsd:simpleType></xsd:element><xsd:element name="MediaServiceKeyPoints" ma:index="11" nillable="true" ma:displayName="KeyPoints" ma:internalName="MediaServiceKeyPoints" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Note"><xsd:maxLength value="255"/></xsd:restriction></xsd:simpleType></xsd:element><xsd:element name="MediaServiceAutoTags" ma:index="14" nillable="true" ma:displayName="Tags" ma:internalName="MediaServiceAutoTags" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element><xsd:element name="MediaServiceOCR" ma:index="15" nillable="true" ma:displayName="Extracted Text" ma:internalName="MediaServiceOCR" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Note"><xsd:maxLength value="255"/></xsd:restriction></xsd:simpleType></xsd:element><xsd:element name="MediaServiceGenerationTime" ma:index="16" nillable="true" ma:displayName="MediaServiceGenerationTime" ma:hidden="true" ma:internalName="MediaServiceGenerationTime" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element><xsd:element name="MediaServiceEventHashCode" ma:index="17" nillable="true" ma:displayName="MediaServiceEventHashCode" ma:hidden="true" ma:internalName="MediaServiceEventHashCode" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element><xsd:element name="MediaServiceDateTaken" ma:index="18" nillable="true" ma:displayName="MediaServiceDateTaken" ma:hidden="true" ma:internalName="MediaServiceDateTaken" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element><xsd:element name="MediaLengthInSeconds" ma:index="19" nillable="true" ma:displayName="MediaLengthInSeconds" ma:hidden="true" ma:internalName="MediaLengthInSeconds" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Unknown"/></xsd:simpleType></xsd:element><xsd:element name="lcf76f155ced4ddcb4097134ff3c332f" ma:index="21" nillable="true" ma:taxonomy="true" ma:internalName="lcf76f155ced4ddcb4097134ff3c332f" ma:taxonomyFieldName="MediaServiceImageTags" ma:displayName="Image Tags" ma:readOnly="false" ma:fieldId="{5cf76f15-5ced-4ddc-b409-7134ff3c332f}" ma:taxonomyMulti="true" ma:sspId="a3ebb39f-d69b-4575-80f5-9912993956e5" ma:termSetId="09814cd3-568e-fe90-9814-8d621ff8fb84" ma:anchorId="fba54fb3-c3e1-fe81-a776-ca4b69148c4d" ma:open="true" ma:isKeyword="false"><xsd:complexType><xsd:sequence><xsd:element ref="pc:Terms" minOccurs="0" maxOccurs="1"/></xsd:sequence></xsd:complexType></xsd:element></xsd:schema><xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="67e6609d-4da6-43f2-b73c-9b749e56ecda" elementFormDefault="qualified"><xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/><xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/><xsd:element name="SharedWithUsers" ma:index="12" nillable="true" ma:displayName="Shared With" ma:internalName="SharedWithUsers" ma:readOnly="true"><xsd:complexType><xsd:complexContent><xsd:extension base="dms:UserMulti"><xsd:sequence><xsd:element name="UserInfo" minOccurs="0" maxOccurs="unbounded"><xsd:complexType><xsd:sequence><xsd:element name="DisplayName" type="xsd:string" minOccurs="0"/><xsd:element name="AccountId" type="dms:UserId" minOccurs="0" nillable="true"/><xsd:element name="AccountType" type="xsd:string" minOccurs="0"/></xsd:sequence></xsd:complexType></xsd:element></xsd:sequence></xsd:extension></xsd:complexContent></xsd:complexType></xsd:element><xsd:element name="SharedWithDetails" ma:index="13" nillable="true" ma:displayName="Shared With Details" ma:internalName="SharedWithDetails" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Note"><xsd:maxLength value="255"/></xsd:restriction></xsd:simpleType></xsd:element><xsd:element name="TaxCatchAll" ma:index="22" nillable="true" ma:displayName="Taxonomy Catch All Column" ma:hidden="true" ma:list="{a6f53592-1060-4e07-aeb7-058c13974b3e}" ma:internalName="TaxCatchAll" ma:showField="CatchAllData" ma:web="67e6609d-4da6-43f2-b73c-9b749e56ecda"><xsd:complexType><xsd:complexContent><xsd:extension base="dms:MultiChoiceLookup"><xsd:sequence><xsd:element name="Value" type="dms:Lookup" maxOccurs="unbounded" minOccurs="0" nillable="true"/></xsd:sequence></xsd:extension></xsd:complexContent></xsd:complexType></xsd:element></xsd:schema><xsd:schema targetNamespace="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" elementFormDefault="qualified" attributeFormDefault="unqualified" blockDefault="#all" xmlns="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:dc="http://purl.org/dc/elements/1.1/" xmlns:dcterms="http://purl.org/dc/terms/" xmlns:odoc="http://schemas.microsoft.com/internal/obd"><xsd:import namespace="http://purl.org/dc/elements/1.1/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd"/><xsd:import namespace="http://purl.org/dc/terms/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd"/><xsd:element name="coreProperties" type="CT_coreProperties"/><xsd:complexType name="CT_coreProperties"><xsd:all><xsd:element ref="dc:creator" minOccurs="0" maxOccurs="1"/><xsd:element ref="dcterms:created" minOccurs="0" maxOccurs="1"/><xsd:element ref="dc:identifier" minOccurs="0" maxOccurs="1"/><xsd:element name="contentType" minOccurs="0" maxOccurs="1" type="xsd:string" ma:index="0" ma:displayName="Content Type"/><xsd:element ref="dc:title" minOccurs="0" maxOccurs="1" ma:index="4" ma:displayName="Title"/><xsd:element ref="dc:subject" minOccurs="0" maxOccurs="1"/><xsd:element ref="dc:description" minOccurs="0" maxOccurs="1"/><xsd:element name="keywords" minOccurs="0" maxOccurs="1" type="xsd:string"/><xsd:element ref="dc:language" minOccurs="0" maxOccurs="1"/><xsd:element name="category" minOccurs="0" maxOccurs="1" type="xsd:string"/><xsd:element name="version" minOccurs="0" maxOccurs="1" type="xsd:string"/><xsd:element name="revision" minOccurs="0" maxOccurs="1" type="xsd:string"><xsd:annotation><xsd:documentation>This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.</xsd:documentation></xsd:annotation></xsd:element><xsd:element name="lastModifiedBy" minOccurs="0" maxOccurs="1" type="xsd:string"/><xsd:element ref="dcterms:modified" minOccurs="0" maxOccurs="1"/><xsd:element name="contentStatus" minOccurs="0" maxOccurs="1" type="xsd:string"/></xsd:all></xsd:complexType></xsd:schema><xs:schema targetNamespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" elementFormDefault="qualified" attributeFormDefault="unqualified" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" xmlns:xs="http://www.w3.org/2001/XMLSchema"><xs:element name="Person"><xs:complexType><xs:sequence><xs:element ref="pc:DisplayName" minOccurs="0"></xs:element><xs:element ref="pc:AccountId" minOccurs="0"></xs:element><xs:element ref="pc:AccountType" minOccurs="0"></xs:element></xs:sequence></xs:complexType></xs:element><xs:element name="DisplayName" type="xs:string"></xs:element><xs:element name="AccountId" type="xs:string"></xs:element><xs:element name="AccountType" type="xs:string"></xs:element><xs:element name="BDCAssociatedEntity"><xs:complexType><xs:sequence><xs:element ref="pc:BDCEntity" minOccurs="0" maxOccurs="unbounded"></xs:element></xs:sequence><xs:attribute ref="pc:EntityNamespace"></xs:attribute><xs:attribute ref="pc:EntityName"></xs:attribute><xs:attribute ref="pc:SystemInstanceName"></xs:attribute><xs:attribute ref="pc:AssociationName"></xs:attribute></xs:complexType></xs:element><xs:attribute name="EntityNamespace" type="xs:string"></xs:attribute><xs:attribute name="EntityName" type="xs:string"></xs:attribute><xs:attribute name="SystemInstanceName" type="xs:string"></xs:attribute><xs:attribute name="AssociationName" type="xs:string"></xs:attribute><xs:element name="BDCEntity"><xs:complexType><xs:sequence><xs:element ref="pc:EntityDisplayName" minOccurs="0"></xs:element><xs:element ref="pc:EntityInstanceReference" minOccurs="0"></xs:element><xs:element ref="pc:EntityId1" minOccurs="0"></xs:element><xs:element ref="pc:EntityId2" minOccurs="0"></xs:element><xs:element ref="pc:EntityId3" minOccurs="0"></xs:element><xs:element ref="pc:EntityId4" minOccurs="0"></xs:element><xs:element ref="pc:EntityId5" minOccurs="0"></xs:element></xs:sequence></xs:complexType></xs:element><xs:element name="EntityDisplayName" type="xs:string"></xs:element><xs:element name="EntityInstanceReference" type="xs:string"></xs:element><xs:element name="EntityId1" type="xs:string"></xs:element><xs:element name="EntityId2" type="xs:string"></xs:element><xs:element name="EntityId3" type="xs:string"></xs:element><xs:element name="EntityId4" type="xs:string"></xs:element><xs:element name="EntityId5" type="xs:string"></xs:element><xs:element name="Terms"><xs:complexType><xs:sequence><xs:element ref="pc:TermInfo" minOccurs="0" maxOccurs="unbounded"></xs:element></xs:sequence></xs:complexType></xs:element><xs:element name="TermInfo"><xs:complexType><xs:sequence><xs:element ref="pc:TermName" minOccurs="0"></xs:element><xs:element ref="pc:TermId" minOccurs="0"></xs:element></xs:sequence></xs:complexType></xs:element><xs:element name="TermName" type="xs:string"></xs:element><xs:element name="TermId" type="xs:string"></xs:element></xs:schema></ct:contentTypeSchema>'
$item3Xml = '<?xml version="1.0" encoding="utf-8"?><p:properties xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"><documentManagement><lcf76f155ced4ddcb4097134ff3c332f xmlns="25558cfc-09cb-4f21-a339-24113151f02d"><Terms xmlns="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"></Terms></lcf76f155ced4ddcb4097134ff3c332f><TaxCatchAll xmlns="67e6609d-4da6-43f2-b73c-9b749e56ecda" xsi:nil="true"/></documentManagement></p:properties>'
$itemProps2Xml = '<?xml version="1.0" encoding="UTF-8" standalone="no"?>
<ds:datastoreItem ds:itemID="{8C0272E7-64E5-46D1-BA9D-607B15B37D08}" xmlns:ds="http://schemas.openxmlformats.org/officeDocument/2006/customXml"/>'

# Replace the three "Document Management" custom XML parts: find the old
# ones by their well-known namespaces, delete them, and add back the
# refreshed XML content (PowerPoint mints new part ids/relationships
# automatically when a custom XML part is added).

$oldSchemaParts = $parts.SelectByNamespace("http://schemas.microsoft.com/office/2006/metadata/contentType")
for ($i = $oldSchemaParts.Count; $i -ge 1; $i--) {
    $oldSchemaParts.Item($i).Delete()
}

$oldPropsParts = $parts.SelectByNamespace("http://schemas.microsoft.com/office/2006/metadata/properties")
for ($i = $oldPropsParts.Count; $i -ge 1; $i--) {
    $oldPropsParts.Item($i).Delete()
}

$null = $parts.Add($item2Xml)
$null = $parts.Add($item3Xml)
$null = $parts.Add($itemProps2Xml)

$p.Save()
